$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6968656182289124
$ws.Range("B1").Value = 1.392574191093445
$ws.Range("C1").Value = 3.753113269805908
$ws.Range("D1").Value = 3.127187967300415
$ws.Range("E1").Value = 1.712586641311646
